$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "Hometown:" -> "Hometown: San Antonio"
$tr.Paragraphs(1).Characters(1, 9).Text = "Hometown: San Antonio"

# Paragraph 2: "Major:" -> "Major: Computer Science"
$tr.Paragraphs(2).Characters(1, 6).Text = "Major: Computer Science"

# Paragraph 3: "Why?" -> "I fell in love with the elective in high school"
$tr.Paragraphs(3).Characters(1, 4).Text = "I fell in love with the elective in high school"
